# Update "想去人数" (number of people interested) figures that changed
# between data pulls, as reflected in the regenerated gh-pages output.
#
# Sheet "展览" (Exhibitions):
#   F2: 984  -> 989
#   F3: 1997 -> 2024
#   F4: 443  -> 445
#
# Sheet "全部类型" (All types) contains the same three records (rows 4-6):
#   F4: 984  -> 989
#   F5: 1997 -> 2024
#   F6: 443  -> 445

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 989
$wsExhibition.Range("F3").Value = 2024
$wsExhibition.Range("F4").Value = 445

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 989
$wsAllTypes.Range("F5").Value = 2024
$wsAllTypes.Range("F6").Value = 445
